# Add "NA" values under the duplicate_image_filename column (column E)
# for data rows 2-21 on the active worksheet ("add the NA's under
# duplicate_image_filename").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}

# Keep the untouched F1 cell blank (the COM round-trip otherwise
# resurrects a stray value there).
$ws.Range("F1").Value = $null
